# Update sheet view / selection state for two worksheets and change the
# workbook's active (visible-on-open) sheet.
#
# Before:
#   - "Frederikse" sheet is the active/tabSelected sheet, scrolled to show
#     its top-left at A1 with selection on G23.
#   - "GrIS" sheet is not active, selection on A9.
#
# After:
#   - "GrIS" sheet becomes the active/tabSelected sheet, selection on A5.
#   - "Frederikse" sheet is no longer active, scrolled so A25 is the
#     top-left visible cell, selection on D44:E44.

$wb = $excel.ActiveWorkbook

$wsFrederikse = $wb.Worksheets.Item("Frederikse")
$wsGrIS = $wb.Worksheets.Item("GrIS")

# Update the Frederikse sheet: scroll position + selection, no longer the
# tab shown when the workbook opens.
$wsFrederikse.Activate()
$wsFrederikse.Range("D44:E44").Select()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1

# Update the GrIS sheet: becomes the active tab, with a new selection.
$wsGrIS.Activate()
$wsGrIS.Range("A5").Select()

$wb.Save()
